# Add a new "24-ago" column (BE) right after the existing "23-ago" column (BD),
# carrying the same per-row values pattern as the rest of the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Cells.Item(1, 57).Value = "24-ago"

# Data values for rows 2-18
$values = @{
    2  = 0
    3  = 11.30610872990993
    4  = 19.613724179130244
    5  = 21.209603009491765
    6  = 0
    7  = 13.847620683361276
    8  = 24.839863034507928
    9  = 11.859902798862157
    10 = 3.7025223923995703
    11 = 13.882607946509887
    12 = 0
    13 = 8.6832817927842729
    14 = 0
    15 = 0
    16 = 4.3194735259753889
    17 = 0
    18 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 57).Value = $values[$row]
}

# Move the selection to match the post-edit workbook state
$ws.Range("BG7").Select()
